$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV & FV")

# --- Row 23: period numbers (same layout as row 7) ---
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 5

# --- Row 24: cash flows (same layout as row 8) ---
$ws.Range("C24").Value = -500
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 120
$ws.Range("F24").Value = 200
$ws.Range("G24").Value = 120
$ws.Range("H24").Value = 120

# --- Row 25: discounted cash flows (PV formulas) ---
$ws.Range("C25").Value = -500
$ws.Range("D25").Formula = '=D24/($C$4+1)^D23'
$ws.Range("E25").Formula = '=E24/($C$4+1)^E23'
$ws.Range("F25").Formula = '=F24/($C$4+1)^F23'
$ws.Range("G25").Formula = '=G24/($C$4+1)^G23'
$ws.Range("H25").Formula = '=H24/($C$4+1)^H23'

# Match the formatting used elsewhere on the sheet for the newly-built
# lookup/database block (copy the number/font formats from the matching
# columns used earlier in the sheet, without disturbing the values/formulas
# just written). Pasted one contiguous range at a time (multi-area
# destinations are unreliable for PasteSpecial in this runtime).
$ws.Range("H8").Copy()
$ws.Range("C23:D23").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("F23:H23").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("C24:D24").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("F24:H24").PasteSpecial(-4122)
$ws.Range("H8").Copy()
$ws.Range("C25:H25").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E8").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 26 / 27 / 31: left-over formatted-but-empty cells (style only) ---
$ws.Range("B8").Copy()
$ws.Range("B26").PasteSpecial(-4122)
$ws.Range("C8").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D8").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("F8").Copy()
$ws.Range("F26:G26").PasteSpecial(-4122)

$ws.Range("F8").Copy()
$ws.Range("D27:H27").PasteSpecial(-4122)

$ws.Range("B20").Copy()
$ws.Range("B31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Sheet view: scroll + selection moved to reflect the new data block ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("B27").Select()
